# Apply Week 17 data updates to "Players Data.xlsx"
# Sheet "Rushing" and Sheet "Receiving" get updated attempt/target/completion counts.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2: J.Herbert
$rushing.Range("E2").Value = 21
$rushing.Range("F2").Value = 13

# Row 3: A.Ekeler
$rushing.Range("C3").Value = 115
$rushing.Range("D3").Value = 62
$rushing.Range("E3").Value = 12
$rushing.Range("F3").Value = 44

# Row 4: J.Jackson
$rushing.Range("C4").Value = 29
$rushing.Range("D4").Value = 22
$rushing.Range("E4").Value = 4
$rushing.Range("F4").Value = 17

# Row 7: G.Nabers
$rushing.Range("C7").Value = 2

# Row 10: J.Guyton
$rushing.Range("C10").Value = 5

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2: A.Ekeler
$receiving.Range("C2").Value = 68
$receiving.Range("D2").Value = 55
$receiving.Range("E2").Value = 5
$receiving.Range("F2").Value = 3

# Row 3: J.Jackson
$receiving.Range("C3").Value = 12
$receiving.Range("D3").Value = 11
$receiving.Range("G3").Value = 5
$receiving.Range("H3").Value = 5

# Row 6: K.Allen
$receiving.Range("C6").Value = 117
$receiving.Range("D6").Value = 79
$receiving.Range("G6").Value = 23
$receiving.Range("H6").Value = 14

# Row 7: M.Williams
$receiving.Range("C7").Value = 70
$receiving.Range("D7").Value = 43
$receiving.Range("E7").Value = 28
$receiving.Range("F7").Value = 13

# Row 8: J.Palmer
$receiving.Range("C8").Value = 25
$receiving.Range("D8").Value = 18
$receiving.Range("G8").Value = 4
$receiving.Range("H8").Value = 3

# Row 9: J.Guyton
$receiving.Range("C9").Value = 26
$receiving.Range("D9").Value = 16
$receiving.Range("G9").Value = 6
$receiving.Range("H9").Value = 3

# Row 13: S.Anderson
$receiving.Range("C13").Value = 13
$receiving.Range("D13").Value = 11

# Row 14: T.McKitty
$receiving.Range("C14").Value = 6
$receiving.Range("D14").Value = 5
